# Add 2022-Q3 data:
#  - Insert a new "2022-Q3" worksheet right after "总计", built by
#    copying the "2022-Q2" sheet (so header/row styles match exactly)
#    and then overwriting the fund-level figures with the new quarter's
#    numbers.
#  - Insert a new summary row for "2022-Q3" at the top of the "总计"
#    sheet's data (row 2), pushing the existing quarters down.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet as a copy of "2022-Q2" ---
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Overwrite the fund figures with the 2022-Q3 values, keeping them as
# text (matching the original sheet's inline-string cell type).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $wsQ3.Cells.Item(2, 4) "10.25"
Set-TextValue $wsQ3.Cells.Item(2, 5) "93.67"
Set-TextValue $wsQ3.Cells.Item(2, 6) "1.10"
Set-TextValue $wsQ3.Cells.Item(2, 7) "0.1128"

Set-TextValue $wsQ3.Cells.Item(3, 4) "1.07"
Set-TextValue $wsQ3.Cells.Item(3, 5) "97.91"
Set-TextValue $wsQ3.Cells.Item(3, 6) "1.16"
Set-TextValue $wsQ3.Cells.Item(3, 7) "0.0124"

# --- 2. Insert a new row into "总计" for 2022-Q3, shifting rows down ---
# (Capture the existing data rows first - note the `()` on Value() is
#  required here, a bare `.Value` property access doesn't invoke the
#  getter in this host.)
$existingRows = @()
for ($r = 2; $r -le 5; $r++) {
    $rowVals = @(
        $wsTotal.Cells.Item($r, 2).Value(),
        $wsTotal.Cells.Item($r, 3).Value(),
        $wsTotal.Cells.Item($r, 4).Value()
    )
    $existingRows += , $rowVals
}

# Push the 4 existing quarters down into rows 3..6, keeping their style.
for ($i = 0; $i -lt $existingRows.Count; $i++) {
    $r = $i + 3
    $rowVals = $existingRows[$i]
    $wsTotal.Cells.Item($r, 2).Value = $rowVals[0]
    $wsTotal.Cells.Item($r, 3).Value = $rowVals[1]
    $wsTotal.Cells.Item($r, 4).Value = $rowVals[2]
    $wsTotal.Cells.Item($r, 1).Value = $i + 1
}

# Make sure the newly-created row 6's index cell (A6) carries the same
# bold/bordered/centered style used by the rest of column A.
$aStyle = $wsTotal.Cells.Item(5, 1)
$a6 = $wsTotal.Cells.Item(6, 1)
$a6.Font.Bold = $true
$a6.HorizontalAlignment = -4108  # xlCenter
$a6.VerticalAlignment = -4160    # xlTop
$a6.Borders.LineStyle = 1

# Fill in the brand-new 2022-Q3 summary row.
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.13
